# "Add files via upload" -- adds a new "Sheet3" at the end of the workbook
# containing the single value "6787-897" in A1, makes it the active/selected
# sheet, and trims the now-unused trailing blank rows (3:4) from Sheet1.

$wb = $excel.ActiveWorkbook

# Sheet1 had two stray blank rows (3 and 4) below its real data -- drop them
# so the sheet's used range shrinks back down to A1:B2.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3:B4").ClearContents()

# Append the new sheet right after Sheet2 (i.e. at the end of the workbook).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"
$ws3.Range("A1").Value = "6787-897"

# Sheet3 becomes the active/selected tab (Sheet2 loses tabSelected).
$ws3.Activate()
$ws3.Select()
